$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from the generic "CHROM" placeholder to the actual
# reference/chromosome used for this sample (te_09-0932).
$ws.Name = "NZ_CP021201.1"

# Column A previously held blank placeholder rows under a generic header.
# Fill in the real header + filtered-region position values.
$ws.Range("A1").Value = "te-09-0932-All"
$ws.Range("A2").Value = "1034828"
$ws.Range("A3").Value = "406316"
$ws.Range("A4").Value = "406337"
$ws.Range("A5").Value = "407417"
$ws.Range("A6").Value = "1437159"
$ws.Range("A7").Value = "1436983"

# New column B: the "-02" sheep sample added alongside the "-All" column.
$ws.Range("B1").Value = "te-09-0932-02"
$ws.Range("B2").Value = "1564433"

# Mirror the author's final UI state: column B selected (whole column),
# as if it had just been added/selected to paste the new sample's data.
[void]$ws.Columns("B").Select()
